$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values in P1/Q1 and copy the existing header format (style) from O1
# so they share the same style index as the rest of the header row instead of
# creating a brand-new style entry.
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Update data rows 2 through 25: swap I/K and M/O column values, and populate
# the two newly-added columns P and Q with value 2.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1

    $ws.Cells.Item($r, 16).Value = 2  # P (new column)
    $ws.Cells.Item($r, 17).Value = 2  # Q (new column)
}
